$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage so numeric-looking strings (e.g. "1.007", "309.31")
# are not auto-converted to numbers by Excel, matching the source data
# which stores these as inline strings, not numeric cells.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.862.79"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "1.841.80"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "309.31"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "0.4704"
$ws.Range("D8").Value = "0.3658"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").Value = "0.07143"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").Value = "0.9249"
$ws.Range("E10").Value = "  +3.86%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "19.53"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07683"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.879.19"
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("D14").Value = "5.287"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "6.391"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "88.19"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "0.000008618"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "26.898.68"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").Value = "5.010"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").Value = "10.59"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").Value = "1.925"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").Value = "151.75"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "18.22"
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").Value = "114.14"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").Value = "4.878"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").Value = "0.08821"
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("D31").Value = "3.211"
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("D32").Value = "1.177"
$ws.Range("E32").Value = "  +6.04%  "
$ws.Range("D33").Value = "0.7462"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").Value = "2.776"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("D35").Value = "4.466"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").Value = "0.01936"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "0.05201"
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("D39").Value = "2.955"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").Value = "0.5198"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").Value = "6.956"
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("D42").Value = "0.1508"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "8.147"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("D44").Value = "10.49"
$ws.Range("E44").Value = "  +5.58%  "
$ws.Range("D45").Value = "0.4695"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "1.005"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").Value = "101.46"
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("D48").Value = "1.594"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").Value = "65.68"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").Value = "0.06033"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "0.8894"
$ws.Range("E51").Value = "  +5.08%  "

# Restore default (unstyled) cell style now that the text values are set,
# so no residual "@" number-format style lingers on these cells.
$dataRange.Style = "Normal"
